$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.188.47"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "3.068.32"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.71"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.61"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.208"
$ws.Range("E9").Value = "  +5.65%  "
$ws.Range("D10").Value = "3.066.89"
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("E13").Value = "  +6.56%  "
$ws.Range("D14").Value = "3.609.68"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.91"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "76.282.86"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000193"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "3.066.82"
$ws.Range("E18").Value = "  +3.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.13"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.72"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").Value = "  +9.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.39"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").Value = "3.240.24"
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.19"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.33"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.84"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "497.67"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.123"
$ws.Range("E36").Value = "  +11.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.67"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.16"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "193.39"
$ws.Range("E40").Value = "  +7.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.377"
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.102"
$ws.Range("E42").Value = "  -9.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.793"
$ws.Range("E44").Value = "  +20.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.10"
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.25"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.594"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.87"
$ws.Range("E51").Value = "  -0.87%  "
